$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple price updates (Price column D) - prefix with apostrophe so the
# values stay text (matching the workbook's inlineStr string storage)
# instead of being auto-converted to numbers.
$ws.Range("D2").Value  = "'246.69"
$ws.Range("D4").Value  = "'5.426"
$ws.Range("D5").Value  = "'0.05782"
$ws.Range("D7").Value  = "'6.332"
$ws.Range("D8").Value  = "'0.8137"
$ws.Range("D9").Value  = "'0.9429"
$ws.Range("D10").Value = "'0.1427"
$ws.Range("D11").Value = "'0.07515"
$ws.Range("D12").Value = "'0.03147"
$ws.Range("D14").Value = "'4.149"
$ws.Range("D15").Value = "'0.09409"
$ws.Range("D16").Value = "'0.001590"
$ws.Range("D17").Value = "'0.04814"
$ws.Range("D18").Value = "'0.0005900"
$ws.Range("D19").Value = "'0.006195"
$ws.Range("D20").Value = "'0.004124"
$ws.Range("D21").Value = "'0.0009977"
$ws.Range("D23").Value = "'3.770"
$ws.Range("D27").Value = "'0.0001290"
$ws.Range("D44").Value = "'0.006536"
$ws.Range("D48").Value = "'0.1423"

# Rows 41-43 got re-ranked: row 41 (BKEXToken) / 42 (CEJI) / 43 (KickToken)
# become KickToken / BKEXToken / CEJI respectively, each with refreshed
# price + rank-label (column E) data. Columns A, F, G are unchanged.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006319"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1075"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003000"
$ws.Range("E43").Value = "42CEJICEJI"
